# Auto-generated: refresh market-price snapshot columns (H-N) across all job sheets
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 178.66667
$ws.Range("I4").Value = 198.5
$ws.Range("K4").Value = 198.5
$ws.Range("M4").Value = -84.5

# Row 53
$ws.Range("H53").Value = 737.7
$ws.Range("I53").Value = 250
$ws.Range("J53").Value = 1469.25
$ws.Range("K53").Value = 250
$ws.Range("L53").Value = 1469.25
$ws.Range("M53").Value = 387
$ws.Range("N53").Value = -2743.25

# Row 112
$ws.Range("H112").Value = 857399
$ws.Range("J112").Value = 1113850.6
$ws.Range("L112").Value = 3341551.8
$ws.Range("N112").Value = -3343767.8

# Row 120
$ws.Range("H120").Value = 83804
$ws.Range("J120").Value = 83804
$ws.Range("L120").Value = 83804
$ws.Range("N120").Value = -93480

# Row 125
$ws.Range("H125").Value = 3898.3333
$ws.Range("J125").Value = 3898.3333
$ws.Range("L125").Value = 35084.9997
$ws.Range("N125").Value = -40004.9997

# Row 132
$ws.Range("H132").Value = 1603.1136
$ws.Range("I132").Value = 1571.3414
$ws.Range("K132").Value = 4714.0242
$ws.Range("M132").Value = -2184.0242

# Row 137
$ws.Range("H137").Value = 11116151
$ws.Range("I137").Value = 6950
$ws.Range("J137").Value = 18522286
$ws.Range("K137").Value = 20850
$ws.Range("L137").Value = 55566858
$ws.Range("M137").Value = -18300
$ws.Range("N137").Value = -55571958

# Row 141
$ws.Range("H141").Value = 7688.2
$ws.Range("I141").Value = 6924.0386
$ws.Range("K141").Value = 20772.1158
$ws.Range("M141").Value = -15592.1158


# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 29414506
$ws.Range("I32").Value = 30569858
$ws.Range("J32").Value = 17860986
$ws.Range("K32").Value = 30569858
$ws.Range("L32").Value = 17860986
$ws.Range("M32").Value = -30569571
$ws.Range("N32").Value = -17861560

# Row 74
$ws.Range("H74").Value = 3085.9546
$ws.Range("I74").Value = 3129.6316
$ws.Range("J74").Value = 2809.3333
$ws.Range("K74").Value = 3129.6316
$ws.Range("L74").Value = 2809.3333
$ws.Range("M74").Value = -2255.6316
$ws.Range("N74").Value = -4557.3333

# Row 77
$ws.Range("H77").Value = 3085.9546
$ws.Range("I77").Value = 3129.6316
$ws.Range("J77").Value = 2809.3333
$ws.Range("K77").Value = 15648.158
$ws.Range("L77").Value = 14046.6665
$ws.Range("M77").Value = -11280.158
$ws.Range("N77").Value = -22782.6665

# Row 110
$ws.Range("H110").Value = 1460.375
$ws.Range("J110").Value = 1954.5454
$ws.Range("L110").Value = 1954.5454
$ws.Range("N110").Value = -6044.5454

# Row 132
$ws.Range("H132").Value = 3591.9143
$ws.Range("I132").Value = 3444.7932
$ws.Range("J132").Value = 4303
$ws.Range("K132").Value = 10334.3796
$ws.Range("L132").Value = 12909
$ws.Range("M132").Value = -7804.3796
$ws.Range("N132").Value = -17969


# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3105.7693
$ws.Range("I20").Value = 2547.125
$ws.Range("K20").Value = 2547.125
$ws.Range("M20").Value = -2300.125

# Row 95
$ws.Range("H95").Value = 132000
$ws.Range("J95").Value = 132000
$ws.Range("L95").Value = 132000
$ws.Range("N95").Value = -137492

# Row 134
$ws.Range("H134").Value = 5955390.5
$ws.Range("I134").Value = 8930884
$ws.Range("K134").Value = 26792652
$ws.Range("M134").Value = -26790117


# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6438.85
$ws.Range("I31").Value = 2900
$ws.Range("K31").Value = 2900
$ws.Range("M31").Value = -2605

# Row 34
$ws.Range("H34").Value = 6438.85
$ws.Range("I34").Value = 2900
$ws.Range("K34").Value = 2900
$ws.Range("M34").Value = -2698

# Row 99
$ws.Range("H99").Value = 2664.6667
$ws.Range("I99").Value = 2664.6667
$ws.Range("K99").Value = 2664.6667
$ws.Range("M99").Value = -1166.6667

# Row 105
$ws.Range("H105").Value = 1657.7273
$ws.Range("I105").Value = 630.9167
$ws.Range("K105").Value = 630.9167
$ws.Range("M105").Value = 1116.0833

# Row 126
$ws.Range("H126").Value = 2664.6667
$ws.Range("I126").Value = 2664.6667
$ws.Range("K126").Value = 7994.000100000001
$ws.Range("M126").Value = -5524.000100000001


# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 3087.7856
$ws.Range("I18").Value = 1923
$ws.Range("J18").Value = 5999.75
$ws.Range("K18").Value = 5769
$ws.Range("L18").Value = 17999.25
$ws.Range("M18").Value = -5600
$ws.Range("N18").Value = -18337.25

# Row 122
$ws.Range("H122").Value = 99.25
$ws.Range("I122").Value = 99.25
$ws.Range("K122").Value = 893.25
$ws.Range("M122").Value = 1556.75

# Row 140
$ws.Range("H140").Value = 2149.4707
$ws.Range("I140").Value = 1769.6
$ws.Range("J140").Value = 4998.5
$ws.Range("K140").Value = 5308.799999999999
$ws.Range("L140").Value = 14995.5
$ws.Range("M140").Value = -128.7999999999993
$ws.Range("N140").Value = -25355.5


# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 107
$ws.Range("H107").Value = 1483.5555
$ws.Range("I107").Value = 1618.6666
$ws.Range("J107").Value = 1416
$ws.Range("K107").Value = 1618.6666
$ws.Range("L107").Value = 1416
$ws.Range("M107").Value = 301.3334
$ws.Range("N107").Value = -5256

# Row 132
$ws.Range("H132").Value = 3109.95
$ws.Range("I132").Value = 2679.1333
$ws.Range("K132").Value = 8037.3999
$ws.Range("M132").Value = -5507.3999


# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4177.9165
$ws.Range("I40").Value = 3737.1765
$ws.Range("J40").Value = 5248.2856
$ws.Range("K40").Value = 3737.1765
$ws.Range("L40").Value = 5248.2856
$ws.Range("M40").Value = -3601.1765
$ws.Range("N40").Value = -5520.2856

# Row 46
$ws.Range("H46").Value = 7625.769
$ws.Range("I46").Value = 2018.8
$ws.Range("J46").Value = 8960.762000000001
$ws.Range("K46").Value = 2018.8
$ws.Range("L46").Value = 8960.762000000001
$ws.Range("M46").Value = -1830.8
$ws.Range("N46").Value = -9336.762000000001

# Row 68
$ws.Range("H68").Value = 3734.375
$ws.Range("I68").Value = 3767.8572
$ws.Range("K68").Value = 3767.8572
$ws.Range("M68").Value = -3018.8572

# Row 71
$ws.Range("H71").Value = 3734.375
$ws.Range("I71").Value = 3767.8572
$ws.Range("K71").Value = 18839.286
$ws.Range("M71").Value = -15095.286

# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 93
$ws.Range("H93").Value = 2512.6
$ws.Range("I93").Value = 2140.75
$ws.Range("K93").Value = 2140.75
$ws.Range("M93").Value = -892.75

# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# Row 95
$ws.Range("H95").Value = 100000
$ws.Range("J95").Value = 100000
$ws.Range("L95").Value = 100000
$ws.Range("N95").Value = -105492

# Row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# Row 98
$ws.Range("H98").Value = 30355
$ws.Range("J98").Value = 30355
$ws.Range("L98").Value = 30355
$ws.Range("N98").Value = -36345

# Row 99
$ws.Range("H99").Value = 48354.6
$ws.Range("J99").Value = 106979
$ws.Range("L99").Value = 106979
$ws.Range("N99").Value = -112969

# Row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# Row 101
$ws.Range("H101").Value = 50837
$ws.Range("J101").Value = 50837
$ws.Range("L101").Value = 50837
$ws.Range("N101").Value = -57327

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# Row 103
$ws.Range("H103").Value = 110000
$ws.Range("J103").Value = 110000
$ws.Range("L103").Value = 110000
$ws.Range("N103").Value = -112344

# Row 104
$ws.Range("H104").Value = 74664.664
$ws.Range("J104").Value = 74664.664
$ws.Range("L104").Value = 74664.664
$ws.Range("N104").Value = -81652.664

# Row 105
$ws.Range("H105").Value = 123000
$ws.Range("J105").Value = 123000
$ws.Range("L105").Value = 123000
$ws.Range("N105").Value = -129988

# Row 106
$ws.Range("H106").Value = 5020000
$ws.Range("J106").Value = 5020000
$ws.Range("L106").Value = 5020000
$ws.Range("N106").Value = -5022524

# Row 122
$ws.Range("H122").Value = 25463.846
$ws.Range("I122").Value = 25670.555
$ws.Range("J122").Value = 24998.75
$ws.Range("K122").Value = 77011.66500000001
$ws.Range("L122").Value = 74996.25
$ws.Range("M122").Value = -74561.66500000001
$ws.Range("N122").Value = -79896.25


# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 854.1111
$ws.Range("I100").Value = 364.83334
$ws.Range("K100").Value = 729.66668
$ws.Range("M100").Value = -188.66668

# Row 107
$ws.Range("H107").Value = 369.7
$ws.Range("I107").Value = 235
$ws.Range("K107").Value = 705
$ws.Range("M107").Value = 1215

# Row 122
$ws.Range("H122").Value = 4243.3076
$ws.Range("I122").Value = 4243.3076
$ws.Range("K122").Value = 12729.9228
$ws.Range("M122").Value = -10279.9228

# Row 136
$ws.Range("H136").Value = 7180078.5
$ws.Range("I136").Value = 918.3158
$ws.Range("J136").Value = 22336082
$ws.Range("K136").Value = 2754.9474
$ws.Range("L136").Value = 67008246
$ws.Range("M136").Value = -204.9474
$ws.Range("N136").Value = -67013346

